# [PV-350][WIP] Replace hard coding of visual height with calculated value
# Update the header row on the "PV-Test-03" sheet to the new column titles.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PV-Test-03")

$ws.Range("A1").Value = "Id"
$ws.Range("C1").Value = "Task Name"
$ws.Range("E1").Value = "Start"
$ws.Range("F1").Value = "Finish"
